$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("TCNA_DL_Ferreres", "Tipo de cambio nominal anual - Dólar Libre (Ferreres = Billetes Argentinos)", "Tipo_Cambio_Arg"),
    @("TCNA_EXPO_Ferreres", "Tipo de cambio nominal anual - Dólar Exportaciones (Ferreres)", "Tipo_Cambio_Arg"),
    @("TCNA_JIC", "Tipo de cambio nominal anual - Iñigo Carrera ", "Tipo_Cambio_Arg"),
    @("TCNA_Arceo", "Tipo de cambio nominal anual - Arceo", "Tipo_Cambio_Arg"),
    @("TCRA_2007_DL_Ferreres", "Tipo de cambio real anual  (2007 = 100)- Dólar Libre (Ferreres = Billetes Argentinos)", "Tipo_Cambio_Arg"),
    @("TCRA_2007_EXPO_Ferreres", "Tipo de cambio real anual  (2007 = 100)- Dólar Exportaciones (Ferreres)", "Tipo_Cambio_Arg"),
    @("TCRA_2007_JIC", "Tipo de cambio real anual  (2007 = 100)- Iñigo Carrera ", "Tipo_Cambio_Arg"),
    @("TCRA_2007_Arceo", "Tipo de cambio real anual  (2007 = 100)- Arceo", "Tipo_Cambio_Arg"),
    @("TCP_IPC_DL_Ferreres", "Tipo de cambio de paridad anual - Metodología IPC - (Base-60-71) - Dólar Libre (Ferreres = Billetes Argentinos)", "Tipo_Cambio_Arg"),
    @("TCP_IPC_EXPO_Ferreres", "Tipo de cambio de paridad anual - Metodología IPC - (Base-60-71) - Dólar Exportaciones (Ferreres)", "Tipo_Cambio_Arg"),
    @("TCP_IPC_JIC", "Tipo de cambio de paridad anual - Metodología IPC - (Base-60-71) - Iñigo Carrera ", "Tipo_Cambio_Arg"),
    @("TCP_IPC_Arceo", "Tipo de cambio de paridad anual - Metodología IPC - (Base-60-71)  - Arceo", "Tipo_Cambio_Arg"),
    @("TCP_IPC_PROD_DL_Ferreres", "Tipo de cambio de paridad anual - Metodología IPC + Productividad - (Base-60-71) - Dólar Libre (Ferreres = Billetes Argentinos)", "Tipo_Cambio_Arg"),
    @("TCP_IPC_PROD_EXPO_Ferreres", "Tipo de cambio de paridad anual - Metodología IPC - (Base-60-71) - Dólar Exportaciones (Ferreres)", "Tipo_Cambio_Arg"),
    @("TCP_IPC_PROD_JIC", "Tipo de cambio de paridad anual - Metodología IPC + Productividad - (Base-60-71) - Iñigo Carrera ", "Tipo_Cambio_Arg"),
    @("TCP_IPC_PROD_Arceo", "Tipo de cambio de paridad anual - Metodología IPC + Productividad - (Base-60-71)  - Arceo", "Tipo_Cambio_Arg"),
    @("GS_IPC_DL_Ferreres", "Grado de sobrevaluación - Metodología IPC - (Base-60-71) - Dólar Libre (Ferreres = Billetes Argentinos)", "Tipo_Cambio_Arg"),
    @("GS_IPC_EXPO_Ferreres", "Grado de sobrevaluación - Metodología IPC - (Base-60-71) - Dólar Exportaciones (Ferreres)", "Tipo_Cambio_Arg"),
    @("GS_IPC_JIC", "Grado de sobrevaluación - Metodología IPC - (Base-60-71) - Iñigo Carrera ", "Tipo_Cambio_Arg"),
    @("GS_IPC_Arceo", "Grado de sobrevaluación - Metodología IPC - (Base-60-71)  - Arceo", "Tipo_Cambio_Arg"),
    @("GS_IPC_PROD_DL_Ferreres", "Grado de sobrevaluación - Metodología IPC + Productividad - (Base-60-71) - Dólar Libre (Ferreres = Billetes Argentinos)", "Tipo_Cambio_Arg"),
    @("GS_IPC_PROD_EXPO_Ferreres", "Grado de sobrevaluación - Metodología IPC - (Base-60-71) - Dólar Exportaciones (Ferreres)", "Tipo_Cambio_Arg"),
    @("GS_IPC_PROD_JIC", "Grado de sobrevaluación - Metodología IPC + Productividad - (Base-60-71) - Iñigo Carrera ", "Tipo_Cambio_Arg"),
    @("GS_IPC_PROD_Arceo", "Grado de sobrevaluación - Metodología IPC + Productividad - (Base-60-71)  - Arceo", "Tipo_Cambio_Arg"),
)

$startRow = 12
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
